# Re-applies an upstream re-shuffle of the per-observation rows in the
# "Artfynd" sheet. Row 5 is untouched; rows 2,3,4,6-17 are a permutation
# of each other's full row content (every column A:AY), keyed off of the
# old "Id" column A value. We reproduce that permutation here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = 51   # column AY
$rows = @(2,3,4,6,7,8,9,10,11,12,13,14,15,16,17)

# Columns that are present-but-blank placeholder cells (inlineStr with no
# text) in every one of the affected rows, both before and after the
# shuffle. Round-tripping them through Range.Value2 turns a "present but
# empty" cell into a genuinely absent one (the COM bridge can't recreate
# that nuance), so we simply never touch them - they are already correct
# as-is since every row in play has them blank on both sides of the edit.
$skipCols = @(9, 46, 51)   # I, AT, AY

# destination row -> source row (content that should end up living at the
# destination row, taken from the *original* workbook layout)
$mapping = @{
    2  = 9
    3  = 10
    4  = 15
    6  = 11
    7  = 3
    8  = 4
    9  = 14
    10 = 8
    11 = 13
    12 = 7
    13 = 16
    14 = 17
    15 = 12
    16 = 2
    17 = 6
}

# 1) Snapshot every source row's cell values (across all used columns)
#    BEFORE any writes happen, since several rows are both a source and a
#    destination in this permutation.
$snapshot = @{}
foreach ($r in $rows) {
    for ($c = 1; $c -le $lastCol; $c++) {
        if ($skipCols -contains $c) { continue }
        $snapshot[[string]$r + "_" + [string]$c] = $ws.Cells.Item($r, $c).Value2
    }
}

# 2) Write the snapshotted values into their new homes. Strings are
#    written with the cell pre-formatted as Text so the write doesn't get
#    silently reinterpreted as a date/number/time by the COM layer (e.g.
#    "2023-08-02" or "00:00").
foreach ($destRow in $rows) {
    $srcRow = $mapping[$destRow]
    for ($c = 1; $c -le $lastCol; $c++) {
        if ($skipCols -contains $c) { continue }
        $val = $snapshot[[string]$srcRow + "_" + [string]$c]
        $cell = $ws.Cells.Item($destRow, $c)
        if ($val -is [string]) {
            $cell.NumberFormat = "@"
        } else {
            $cell.NumberFormat = "General"
        }
        $cell.Value2 = $val
    }
}
